$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (updated cryptocurrency price/volume data)
$updates = @{
    'D2' = '294.56'
    'E2' = '0.49%'
    'D3' = '31.05'
    'E3' = '0.35%'
    'D4' = '4.938'
    'E4' = '1.60%'
    'D5' = '0.07346'
    'E5' = '1.08%'
    'D6' = '2.308'
    'E6' = '31.53%'
    'D7' = '7.721'
    'E7' = '0.60%'
    'D8' = '3.740'
    'E8' = '-0.59%'
    'D9' = '0.9030'
    'E9' = '0.10%'
    'D10' = '0.1683'
    'E10' = '1.81%'
    'D11' = '0.07992'
    'E11' = '5.34%'
    'D12' = '0.08187'
    'E12' = '0.12%'
    'D13' = '0.03102'
    'E13' = '2.16%'
    'E14' = '0.82%'
    'D15' = '0.001499'
    'D16' = '0.005814'
    'E16' = '2.83%'
    'D17' = '3.481'
    'E17' = '0.66%'
    'E18' = '-1.52%'
    'E19' = '1.07%'
    'E20' = '-0.23%'
    'D21' = '3.989'
    'E21' = '-8.56%'
    'E22' = '4.60%'
    'D23' = '0.04527'
    'E23' = '0.78%'
    'E24' = '-0.52%'
    'D25' = '0.004658'
    'E25' = '15.30%'
    'E26' = '3.64%'
    'D39' = '0.01609'
    'E39' = '-2.49%'
    'D40' = '0.04453'
    'E40' = '1.90%'
    'D41' = '0.007362'
    'E41' = '-0.63%'
    'D42' = '0.1327'
    'E42' = '0.56%'
    'D43' = '0.008577'
    'E44' = '-2.13%'
    'D45' = '0.009501'
    'E45' = '-7.38%'
    'D46' = '0.00005904'
    'E46' = '4.14%'
    'D47' = '0.00000000748'
    'E47' = '-0.45%'
    'E48' = '3.10%'
    'E49' = '20.39%'
    'E50' = '-0.45%'
    'E51' = '-0.45%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
